$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 71, shifting the existing rows 71..185 down to 72..186.
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new record. It shares all the
# "static" field values with the (now shifted-down) row that used to be row 71
# (A=3, B=Femacal de La Calera, C=Coquimbo, ..., I=Primera, K=14000, L=15000,
# N=$/saco 25 kilos, O=Provincia de Limari, Q=25, R=Hortaliza) except for the
# date (D), volume (J), weighted average price (M) and $/Kg price (P), which
# differ for this new entry.
$ws.Range("A71").Value2 = 3
$ws.Range("B71").Value2 = "Femacal de La Calera"
$ws.Range("C71").Value2 = "Coquimbo"
$ws.Range("D71").Value2 = 44797
$ws.Range("E71").Value2 = 5
$ws.Range("F71").Value2 = 100112026
$ws.Range("G71").Value2 = "Haba"
$ws.Range("H71").Value2 = "Sin especificar"
$ws.Range("I71").Value2 = "Primera"
$ws.Range("J71").Value2 = 95
$ws.Range("K71").Value2 = 14000
$ws.Range("L71").Value2 = 15000
$ws.Range("M71").Value2 = 14526
$ws.Range("N71").Value2 = "`$/saco 25 kilos"
$ws.Range("O71").Value2 = "Provincia de Limar" + [char]0x00ED
$ws.Range("P71").Value2 = 581
$ws.Range("Q71").Value2 = 25
$ws.Range("R71").Value2 = "Hortaliza"

# Give the new row's date cell (D71) the same number format as the rest of
# the date column (style index 2 in the original workbook formats D as a
# date/time value).
$ws.Range("D71").NumberFormat = $ws.Range("D72").NumberFormat
